$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "REACH Tuition Course Fees 2021.pdf"
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/sismic/reach/Reach_Tuition_Course_Fees_2021_v1.0.pdf"

$ws.Range("A3").Value = "REACH Offshore Q4 Promotions Region 2.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/offshore/sismic/reach/Reach-Offshore-SISMIC-Q4-Promotions-1OCT-31DEC21_VOL-1.2.pdf"

$ws.Range("B3").Select() | Out-Null
